$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update event names in column A (rows 2-9) to reflect the new
# Daylight Savings Time schedule / naming convention.
$ws.Range("A2").Value = "Defend the Vault 1"
$ws.Range("A3").Value = "Defend the Vault 2"
$ws.Range("A4").Value = "Battlegrounds 1"
$ws.Range("A5").Value = "Battlegrounds 2"
$ws.Range("A6").Value = "Battlegrounds 3"
$ws.Range("A7").Value = "Battlegrounds 4"
$ws.Range("A8").Value = "Corvus Expedition"
$ws.Range("A9").Value = "Rite of Exile"

# Move/restore the active selection to A9.
$ws.Range("A9").Select()
